$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 260.1928506666666
$ws.Range("H2").Value = 780.5785519999999
$ws.Range("I2").Value = 0.8191301249666345
$ws.Range("J2").Value = 0.8191301249666346
$ws.Range("M2").Value = 44.40220133333333
$ws.Range("N2").Value = 133.206604
$ws.Range("O2").Value = 0.9893265572082102
$ws.Range("P2").Value = 0.9893265572082101
$ws.Range("Q2").Value = 11553.13534079526
$ws.Range("R2").Value = 103978.2180671574
$ws.Range("S2").Value = 0.8103871864387715
$ws.Range("T2").Value = 0.8103871864387715
$ws.Range("G3").Value = 260.1928506666666
$ws.Range("H3").Value = 780.5785519999999
$ws.Range("I3").Value = 0.8191301249666345
$ws.Range("J3").Value = 0.8191301249666346
$ws.Range("M3").Value = 0.401961
$ws.Range("N3").Value = 1.205883
$ws.Range("O3").Value = 0.008956103083191794
$ws.Range("P3").Value = 0.008956103083191792
$ws.Range("Q3").Value = 104.587378446824
$ws.Range("R3").Value = 941.2864060214159
$ws.Range("S3").Value = 0.007336213837748955
$ws.Range("T3").Value = 0.007336213837748954
$ws.Range("G4").Value = 260.1928506666666
$ws.Range("H4").Value = 780.5785519999999
$ws.Range("I4").Value = 0.8191301249666345
$ws.Range("J4").Value = 0.8191301249666346
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07707633333333333
$ws.Range("N4").Value = 0.231229
$ws.Range("O4").Value = 0.00171733970859806
$ws.Range("P4").Value = 0.00171733970859806
$ws.Range("Q4").Value = 20.05471088893422
$ws.Range("R4").Value = 180.492398000408
$ws.Range("S4").Value = 0.001406724690114093
$ws.Range("T4").Value = 0.001406724690114093
$ws.Range("I5").Value = 0.1763970508574364
$ws.Range("J5").Value = 0.1763970508574364
$ws.Range("M5").Value = 44.40220133333333
$ws.Range("N5").Value = 133.206604
$ws.Range("O5").Value = 0.9893265572082102
$ws.Range("P5").Value = 0.9893265572082101
$ws.Range("Q5").Value = 2487.930720843797
$ws.Range("R5").Value = 22391.37648759417
$ws.Range("S5").Value = 0.1745142870264692
$ws.Range("T5").Value = 0.1745142870264691
$ws.Range("I6").Value = 0.1763970508574364
$ws.Range("J6").Value = 0.1763970508574364
$ws.Range("M6").Value = 0.401961
$ws.Range("N6").Value = 1.205883
$ws.Range("O6").Value = 0.008956103083191794
$ws.Range("P6").Value = 0.008956103083191792
$ws.Range("Q6").Value = 22.522557225791
$ws.Range("R6").Value = 202.703015032119
$ws.Range("S6").Value = 0.001579830171050226
$ws.Range("T6").Value = 0.001579830171050226
$ws.Range("I7").Value = 0.1763970508574364
$ws.Range("J7").Value = 0.1763970508574364
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07707633333333333
$ws.Range("N7").Value = 0.231229
$ws.Range("O7").Value = 0.00171733970859806
$ws.Range("P7").Value = 0.00171733970859806
$ws.Range("Q7").Value = 4.318717806588556
$ws.Range("R7").Value = 38.868460259297
$ws.Range("S7").Value = 0.0003029336599170671
$ws.Range("T7").Value = 0.0003029336599170671
$ws.Range("G8").Value = 0.6130636666666667
$ws.Range("H8").Value = 1.839191
$ws.Range("I8").Value = 0.00193002581201784
$ws.Range("J8").Value = 0.00193002581201784
$ws.Range("M8").Value = 44.40220133333333
$ws.Range("N8").Value = 133.206604
$ws.Range("O8").Value = 0.9893265572082102
$ws.Range("P8").Value = 0.9893265572082101
$ws.Range("Q8").Value = 27.22137635748489
$ws.Range("R8").Value = 244.992387217364
$ws.Range("S8").Value = 0.00190942579192659
$ws.Range("T8").Value = 0.00190942579192659
$ws.Range("G9").Value = 0.6130636666666667
$ws.Range("H9").Value = 1.839191
$ws.Range("I9").Value = 0.00193002581201784
$ws.Range("J9").Value = 0.00193002581201784
$ws.Range("M9").Value = 0.401961
$ws.Range("N9").Value = 1.205883
$ws.Range("O9").Value = 0.008956103083191794
$ws.Range("P9").Value = 0.008956103083191792
$ws.Range("Q9").Value = 0.246427684517
$ws.Range("R9").Value = 2.217849160653
$ws.Range("S9").Value = 0.00001728551012565272
$ws.Range("T9").Value = 0.00001728551012565272
$ws.Range("G10").Value = 0.6130636666666667
$ws.Range("H10").Value = 1.839191
$ws.Range("I10").Value = 0.00193002581201784
$ws.Range("J10").Value = 0.00193002581201784
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.07707633333333333
$ws.Range("N10").Value = 0.231229
$ws.Range("O10").Value = 0.00171733970859806
$ws.Range("P10").Value = 0.00171733970859806
$ws.Range("Q10").Value = 0.04725269952655555
$ws.Range("R10").Value = 0.425274295739
$ws.Range("S10").Value = 0.000003314509965597452
$ws.Range("T10").Value = 0.000003314509965597452
$ws.Range("G11").Value = 0.5097843333333333
$ws.Range("H11").Value = 1.529353
$ws.Range("I11").Value = 0.001604885390199778
$ws.Range("J11").Value = 0.001604885390199778
$ws.Range("M11").Value = 44.40220133333333
$ws.Range("N11").Value = 133.206604
$ws.Range("O11").Value = 0.9893265572082102
$ws.Range("P11").Value = 0.9893265572082101
$ws.Range("Q11").Value = 22.63554660524577
$ws.Range("R11").Value = 203.719919447212
$ws.Range("S11").Value = 0.001587755737800101
$ws.Range("T11").Value = 0.001587755737800101
$ws.Range("G12").Value = 0.5097843333333333
$ws.Range("H12").Value = 1.529353
$ws.Range("I12").Value = 0.001604885390199778
$ws.Range("J12").Value = 0.001604885390199778
$ws.Range("M12").Value = 0.401961
$ws.Range("N12").Value = 1.205883
$ws.Range("O12").Value = 0.008956103083191794
$ws.Range("P12").Value = 0.008956103083191792
$ws.Range("Q12").Value = 0.204913420411
$ws.Range("R12").Value = 1.844220783699
$ws.Range("S12").Value = 0.00001437351899133769
$ws.Range("T12").Value = 0.00001437351899133769
$ws.Range("G13").Value = 0.5097843333333333
$ws.Range("H13").Value = 1.529353
$ws.Range("I13").Value = 0.001604885390199778
$ws.Range("J13").Value = 0.001604885390199778
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.07707633333333333
$ws.Range("N13").Value = 0.231229
$ws.Range("O13").Value = 0.00171733970859806
$ws.Range("P13").Value = 0.00171733970859806
$ws.Range("Q13").Value = 0.03929230720411111
$ws.Range("R13").Value = 0.353630764837
$ws.Range("S13").Value = 0.000002756133408338971
$ws.Range("T13").Value = 0.000002756133408338971
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2979236666666666
$ws.Range("H14").Value = 0.893771
$ws.Range("I14").Value = 0.0009379129737112659
$ws.Range("J14").Value = 0.000937912973711266
$ws.Range("M14").Value = 44.40220133333333
$ws.Range("N14").Value = 133.206604
$ws.Range("O14").Value = 0.9893265572082102
$ws.Range("P14").Value = 0.9893265572082101
$ws.Range("Q14").Value = 13.22846662929822
$ws.Range("R14").Value = 119.056199663684
$ws.Range("S14").Value = 0.0009279022132426812
$ws.Range("T14").Value = 0.0009279022132426812
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2979236666666666
$ws.Range("H15").Value = 0.893771
$ws.Range("I15").Value = 0.0009379129737112659
$ws.Range("J15").Value = 0.000937912973711266
$ws.Range("M15").Value = 0.401961
$ws.Range("N15").Value = 1.205883
$ws.Range("O15").Value = 0.008956103083191794
$ws.Range("P15").Value = 0.008956103083191792
$ws.Range("Q15").Value = 0.119753694977
$ws.Range("R15").Value = 1.077783254793
$ws.Range("S15").Value = 0.000008400045275621052
$ws.Range("T15").Value = 0.000008400045275621051
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2979236666666666
$ws.Range("H16").Value = 0.893771
$ws.Range("I16").Value = 0.0009379129737112659
$ws.Range("J16").Value = 0.000937912973711266
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07707633333333333
$ws.Range("N16").Value = 0.231229
$ws.Range("O16").Value = 0.00171733970859806
$ws.Range("P16").Value = 0.00171733970859806
$ws.Range("Q16").Value = 0.02296286383988889
$ws.Range("R16").Value = 0.206665774559
$ws.Range("S16").Value = 0.000001610715192963646
$ws.Range("T16").Value = 0.000001610715192963646
